$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the existing hyperlink before the row-shift so we can re-add
# --- it afterwards at its new location (row insert does not relocate the
# --- Hyperlinks collection's stored Range automatically).
$ws.Range("A40").Hyperlinks.Delete()

# --- Insert 6 rows just above the old "Sector Distribution Details" block
# --- (old row 21) so everything from there on shifts down by 6, matching
# --- the target layout (old row 21 -> new row 27, etc).
$ws.Rows("21:26").Insert()

# --- New MSME size-classification table (rows 20-24) ---
$ws.Range("B20").Value = "Number of employees"
$ws.Range("C20").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D20").Value = "Turnover (local currency, unless noted otherwise)"

$ws.Range("A21").Value = "Micro"
$ws.Range("B21").Value = "'"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'"
$ws.Range("D21").Style = "Normal"

$ws.Range("A22").Value = "Small"
$ws.Range("B22").Value = "'"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'"
$ws.Range("D22").Style = "Normal"

$ws.Range("A23").Value = "Medium"
$ws.Range("B23").Value = "Mfg <100 <br/>NonMfg <50"
$ws.Range("C23").Value = "'"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'"
$ws.Range("D23").Style = "Normal"

$ws.Range("A24").Value = "Large"
$ws.Range("B24").Value = "Mfg >=100<br/>NonMfg >=50"
$ws.Range("C24").Value = "'"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'"
$ws.Range("D24").Style = "Normal"

# --- Re-create the hyperlink at its new row (old A40 -> new A46) ---
$ws.Hyperlinks.Add($ws.Range("A46"), "http://www.success.tid.gov.hk/english/lin_sup_org/gov_dep/service_detail_6863.html") | Out-Null
